$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the unit label from "mb" to "pb" for every data row.
$ws.Range("K2:K54").Value = "pb"

# Match the cursor/selection update seen in the diff (K2 -> K2:K54).
$ws.Range("K2:K54").Select() | Out-Null
